$d = $word.ActiveDocument

# Hunk 1: paragraph after the title CENTRO DE INTENDENCIA... gets an empty run with
# rFonts Carlito + color FF0000 (matches paragraph mark formatting already set in pPr)
$d.Paragraphs.Item(9).Range.Font.Name = "Carlito"
$d.Paragraphs.Item(9).Range.Font.Color = 255

# Hunk 2: merge the "contratual; " and "e" runs into a single run "contratual; e"
$rng = $d.Content
[void]$rng.Find.Execute("contratual; e", $true, $false, $false, $false, $false, $true, 1, $false, "contratual; e", 2)

# Hunk 3: empty run inside List Paragraph gets rFonts Carlito + sz/szCs 24
$d.Paragraphs.Item(33).Range.Font.Name = "Carlito"
$d.Paragraphs.Item(33).Range.Font.Size = 12

# Hunk 4: add color=000000 across most of the "Art. 3º Esta Portaria..." sentence
# (everything from the space after "3º" through the end of the paragraph),
# plus the next 3 empty paragraphs and the {{ordenador_despesas}} paragraph.
$d.Range(1617, 1752).Font.Color = 0
$d.Paragraphs.Item(36).Range.Font.Color = 0
$d.Paragraphs.Item(37).Range.Font.Color = 0
$d.Paragraphs.Item(38).Range.Font.Color = 0
$d.Paragraphs.Item(39).Range.Font.Color = 0
